$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column must stay text (values like '1.00' / '51.548.08' are
# display strings, not numbers) - force Text format before assigning so
# COM doesn't auto-coerce them to numeric doubles.

# --- Rows with both Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.632.67"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.985.97"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "383.29"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.64"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.15"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.453.19"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.35"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.60"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.982.10"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +6.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.524.41"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.28"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.85"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.25"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.77"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.28"
$ws.Range("E26").Value = "  +15.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").Value = "  +16.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("E28").Value = "  +13.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.90"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.72"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.05"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0452"
$ws.Range("E37").Value = "  +6.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.99"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.81"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.277"
$ws.Range("E46").Value = "  +11.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.034.22"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0333"
$ws.Range("E51").Value = "  +4.76%  "

# --- Rows with only Volume(1h) (E) updates ---
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("E47").Value = "  -2.28%  "

# --- Rows where Coin name/link/price/volume swapped with adjacent row ---
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.02"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.34"
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").Value = "  +2.71%  "
